# Update column A values in the "localizacao" sheet so the id sequence
# becomes a single continuous series (1..121) across rows 2..122 instead
# of resetting to 0 partway through (at row 62).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 122; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}
